$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing column-A numbering style (used on rows above, e.g. A17)
# onto the filler rows, then fill in sequential ordinal numbers 1..13 for
# rows 20-32 (the filler sentences section).
$ws.Range("A17").Copy()
for ($i = 20; $i -le 32; $i++) {
    $ws.Range("A$i").PasteSpecial(-4122)
    $ws.Cells.Item($i, 1).Value = $i - 19
}

# Update sheet view state: zoom level and selection position
$ws.Range("D32").Select()
$excel.ActiveWindow.Zoom = 114
